$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 15.2749114326304
$ws.Range("C2").Value = 11.14585704520611
$ws.Range("D2").Value = 4.123700783610976
$ws.Range("F2").Value = 18.15325025705543
$ws.Range("G2").Value = 3.594606405955274
$ws.Range("O2").Value = 16.39943604742586

$ws.Range("B3").Value = 14.4444409033573
$ws.Range("C3").Value = 10.54437033021931
$ws.Range("D3").Value = 4.021035508234477
$ws.Range("F3").Value = 18.29609268713367
$ws.Range("G3").Value = 3.596567197462766
$ws.Range("O3").Value = 16.5639995779716

$ws.Range("B4").Value = 13.9085644219173
$ws.Range("C4").Value = 10.15589706423859
$ws.Range("D4").Value = 3.956390498626568
$ws.Range("F4").Value = 18.39294738512732
$ws.Range("G4").Value = 3.597833377876968
$ws.Range("O4").Value = 16.67170413749342

$ws.Range("B5").Value = 13.68382637728384
$ws.Range("C5").Value = 9.99288380487884
$ws.Range("D5").Value = 3.92967103798756
$ws.Range("F5").Value = 18.43469701520658
$ws.Range("G5").Value = 3.598365058078595
$ws.Range("O5").Value = 16.71726394132858

$ws.Range("B6").Value = 13.64612991547485
$ws.Range("C6").Value = 9.965534944671111
$ws.Range("D6").Value = 3.925212438270397
$ws.Range("F6").Value = 18.44176666303982
$ws.Range("G6").Value = 3.59845429293829
$ws.Range("O6").Value = 16.72492976866917

$ws.Range("B7").Value = 13.90555905108783
$ws.Range("C7").Value = 10.1537175075522
$ws.Range("D7").Value = 3.956031634754166
$ws.Range("F7").Value = 18.39350122787943
$ws.Range("G7").Value = 3.597840484662471
$ws.Range("O7").Value = 16.67231182150252

$ws.Range("B8").Value = 14.99404295092689
$ws.Range("C8").Value = 10.94249999418525
$ws.Range("D8").Value = 4.088651965990889
$ws.Range("F8").Value = 18.2005914817285
$ws.Range("G8").Value = 3.59526959718918
$ws.Range("O8").Value = 16.45479128133134

$ws.Range("B9").Value = 16.91684527368789
$ws.Range("C9").Value = 12.33349412972024
$ws.Range("D9").Value = 4.334811370813068
$ws.Range("F9").Value = 17.89577614236179
$ws.Range("G9").Value = 3.590719745717805
$ws.Range("O9").Value = 16.08136564808674

$ws.Range("B10").Value = 18.19460297371116
$ws.Range("C10").Value = 13.25658541941274
$ws.Range("D10").Value = 4.505779693481492
$ws.Range("F10").Value = 17.71782551596267
$ws.Range("G10").Value = 3.58767353331738
$ws.Range("O10").Value = 15.83978988951832

$ws.Range("B11").Value = 18.74575696631755
$ws.Range("C11").Value = 13.6545058125621
$ws.Range("D11").Value = 4.581151280283613
$ws.Range("F11").Value = 17.64711153724361
$ws.Range("G11").Value = 3.586351454641131
$ws.Range("O11").Value = 15.73709947297051

$ws.Range("B12").Value = 18.95008575243374
$ws.Range("C12").Value = 13.80199192777326
$ws.Range("D12").Value = 4.609327768975282
$ws.Range("F12").Value = 17.6218263523392
$ws.Range("G12").Value = 3.585859921010195
$ws.Range("O12").Value = 15.69925742479541

$ws.Range("B13").Value = 18.9062754531393
$ws.Range("C13").Value = 13.77037083186528
$ws.Range("D13").Value = 4.60327598873432
$ws.Range("F13").Value = 17.62720525868408
$ws.Range("G13").Value = 3.58596537716894
$ws.Range("O13").Value = 15.70736078577619

$ws.Range("B14").Value = 18.76265528149552
$ws.Range("C14").Value = 13.66670383799459
$ws.Range("D14").Value = 4.583476799902753
$ws.Range("F14").Value = 17.6450012779704
$ws.Range("G14").Value = 3.586310833575963
$ws.Range("O14").Value = 15.73396519854062

$ws.Range("B15").Value = 18.67411186825888
$ws.Range("C15").Value = 13.60278749890171
$ws.Range("D15").Value = 4.571301127903596
$ws.Range("F15").Value = 17.6560968478
$ws.Range("G15").Value = 3.586523620509275
$ws.Range("O15").Value = 15.75039747063945

$ws.Range("B16").Value = 18.15797322726622
$ws.Range("C16").Value = 13.23013458929281
$ws.Range("D16").Value = 4.500803918869119
$ws.Range("F16").Value = 17.72265460305593
$ws.Range("G16").Value = 3.587761211373546
$ws.Range("O16").Value = 15.84664664055592

$ws.Range("B17").Value = 17.83358712572949
$ws.Range("C17").Value = 12.99586256254453
$ws.Range("D17").Value = 4.456926260986593
$ws.Range("F17").Value = 17.7661223052707
$ws.Range("G17").Value = 3.588536705550185
$ws.Range("O17").Value = 15.9075434799837

$ws.Range("B18").Value = 17.64417748596939
$ws.Range("C18").Value = 12.85904634163408
$ws.Range("D18").Value = 4.431463888830438
$ws.Range("F18").Value = 17.79208608964764
$ws.Range("G18").Value = 3.588988743819026
$ws.Range("O18").Value = 15.94324722329303

$ws.Range("B19").Value = 17.57956215493369
$ws.Range("C19").Value = 12.81236852977881
$ws.Range("D19").Value = 4.42280473115788
$ws.Range("F19").Value = 17.80104162293125
$ws.Range("G19").Value = 3.589142827157549
$ws.Range("O19").Value = 15.95545199822014

$ws.Range("B20").Value = 17.8684120308822
$ws.Range("C20").Value = 13.0210156442725
$ws.Range("D20").Value = 4.461620557443833
$ws.Range("F20").Value = 17.7613953475707
$ws.Range("G20").Value = 3.58845353278793
$ws.Range("O20").Value = 15.90099072771027

$ws.Range("B21").Value = 18.80495926944439
$ws.Range("C21").Value = 13.69724036433358
$ws.Range("D21").Value = 4.589302357481544
$ws.Range("F21").Value = 17.6397334829076
$ws.Range("G21").Value = 3.586209117801643
$ws.Range("O21").Value = 15.72612241647234

$ws.Range("B22").Value = 19.39149363363905
$ws.Range("C22").Value = 14.12054328829089
$ws.Range("D22").Value = 4.67061388713149
$ws.Range("F22").Value = 17.56893005801706
$ws.Range("G22").Value = 3.58479533560021
$ws.Range("O22").Value = 15.61792883023891

$ws.Range("B23").Value = 19.08080153168478
$ws.Range("C23").Value = 13.89633420437526
$ws.Range("D23").Value = 4.627417845147614
$ws.Range("F23").Value = 17.60591554363234
$ws.Range("G23").Value = 3.585545056448855
$ws.Range("O23").Value = 15.67511329355796

$ws.Range("B24").Value = 17.85267676910838
$ws.Range("C24").Value = 13.00965056822629
$ws.Range("D24").Value = 4.459499001224774
$ws.Range("F24").Value = 17.76352937393098
$ws.Range("G24").Value = 3.588491115888616
$ws.Range("O24").Value = 15.90395106792012

$ws.Range("B25").Value = 16.42003072636713
$ws.Range("C25").Value = 11.97433169963323
$ws.Range("D25").Value = 4.269860217175219
$ws.Range("F25").Value = 17.97024533083893
$ws.Range("G25").Value = 3.591898294306101
$ws.Range("O25").Value = 16.17666015713252

